$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text / add new "Tahun Lulus" column (import data & add instansi) ---
$ws.Range("C1").Value = "Jenis Kelamin (Laki-laki/Perempuan)"
$ws.Range("E1").Value = "Tanggal Lulus (dd/M)"
$ws.Range("F1").Value = "Tahun Lulus"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664
$ws.Columns.Item(5).ColumnWidth = 19.333333333333336
$ws.Columns.Item(6).ColumnWidth = 14.333333333333332

# --- Header styling: yellow fill + thin border around all header cells ---
$a1 = $ws.Range("A1")
$a1.Interior.Color = 65535
$a1.Borders.LineStyle = 1
$a1.Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection ---
$ws.Range("D3").Select()
